# Commit message: "Added spaces around --"
#
# Two textual changes in the final paragraph of the body:
#   1. "rate–limiting" -> "rate – limiting"  (spaces added around the en dash)
#   2. A trailing space appended after "affects this factor." (end of document body)

$d = $word.ActiveDocument

$enDash = [char]0x2013

# 1) Add spaces around the en dash in "rate-limiting"
$d.Content.Find.Execute(
    "rate" + $enDash + "limiting",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "rate " + $enDash + " limiting",
    2
)

# 2) Append a trailing space after the very last sentence of the letter
$d.Content.Find.Execute(
    "affects this factor.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "affects this factor. ",
    2
)
